$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Data" sheet: add two new most-recent years (2023, 2022) at the top of
#    the time series and refresh every historical value with the revised
#    figures from the updated source (act tablas web jul25).
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

# Insert two blank rows right after the header row (row 1), pushing the
# existing series down by two rows.
$wsData.Rows.Item(2).Insert()
$wsData.Rows.Item(2).Insert()

$years = @("2023","2022","2021","2020","2019","2018","2017","2016","2015","2014","2013","2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002","2001","2000","1999","1998","1997","1996","1995","1994","1993","1992","1991","1990")
$vals  = @(1417, 1346, 1307, 1375, 1374, 1370, 1336, 1235, 1163, 1059, 1045, 981.4, 928.5, 918.4, 849.6, 810.5, 766.7, 761.1, 429.5, 428.9, 517.6, 605.3, 621.3, 628, 610.7, 948.7, 928.1, 845.9, 722.5, 709.4, 585.1, 482.4, 388.8, 322.8)

for ($i = 0; $i -lt $years.Length; $i++) {
  $r = $i + 2
  # Years are stored as text in column A (like the rest of the sheet), so
  # force a text number format before assigning, otherwise a pure-digit
  # string would be auto-coerced to a number.
  $wsData.Cells.Item($r, 1).NumberFormat = "@"
  $wsData.Cells.Item($r, 1).Value = $years[$i]
  $wsData.Cells.Item($r, 2).Value = $vals[$i]
}

# ---------------------------------------------------------------------------
# 2) "Metadata" sheet: rename the indicator (seguridad social -> proteccion
#    social) and record the new update date.
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Range("B2").Value = "Gasto público anual per cápita en protección social (en USD constantes de 2010)"
$wsMeta.Range("B4").Value = "Gasto público anual per cápita en protección social"

# Insert a new "actualizacion" row right before the "cita" row.
$wsMeta.Rows.Item(9).Insert()
$wsMeta.Range("A9").Value = "actualizacion"
$wsMeta.Range("B9").Value = "Julio 2025"
